$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that are no longer present in the updated data set.
# (Delete from the bottom up so earlier row numbers stay valid.)
$ws.Rows(9).Delete() | Out-Null   # 2023-06-30_00:00:00.000_Progress_Notes_91265
$ws.Rows(5).Delete() | Out-Null   # 2021-12-23_00:00:00.000_Progress_Notes_91107 (6/2021 entry)
$ws.Rows(4).Delete() | Out-Null   # 2021-12-15_00:00:00.000_Progress_Notes_91091

# NOTE: dates are prefixed with a leading apostrophe so Excel stores them
# as plain text (e.g. "2024-02-22") instead of auto-converting them to a
# date serial number, matching the original workbook's text-based dates.

# Row 2 (2024-04-17_00:00:00.000_IMTX_Conference_Note_91586): normalize date, expand evidence sentence
$ws.Range("E2").Value = "'2024-02-22"
$ws.Range("F2").Value = "['Latest Reference Range & Units 02/22/24 09:20   Kappa Free Light Chain 0.76 - 6.83 mg/dL 56.21 (H)   Lambda Free Light Chain 0.68 - 4.58 mg/dL <3.08 (L)   Kappa/Lambda FLC Ratio 0.66 - 2.22  >115.57 (H)']"

# Row 3 (2024-02-08_00:00:00.000_Progress_Notes_91427): fill in missing date, expand evidence sentence
$ws.Range("E3").Value = "'2024-02-08"
$ws.Range("F3").Value = "['Free Light Chains   Result Value Ref Range    Kappa Free Light Chain 203.94 (H) 0.76 - 6.83 mg/dL    Lambda Free Light Chain <0.15 (L) 0.68 - 4.58 mg/dL    Kappa/Lambda FLC Ratio >1456.71 (H) 0.66 - 2.22']"

# Row 4 (now 2024-06-25_00:00:00.000_Progress_Notes_91596, kappa 0.08): normalize date, expand evidence sentence
$ws.Range("E4").Value = "'2024-06-13"
$ws.Range("F4").Value = "['Free Light Chains   Result Value Ref Range    Kappa Free Light Chain 0.08 (L) 0.76 - 6.83 mg/dL    Lambda Free Light Chain <0.15 (L) 0.68 - 4.58 mg/dL    Kappa/Lambda FLC Ratio >0.57 (L) 0.66 - 2.22']"

# Row 5 (now 2021-12-23_00:00:00.000_Progress_Notes_91107, kappa 1.24): normalize date
$ws.Range("E5").Value = "'2021-12-15"

# Row 6 (now 2022-01-18_00:00:00.000_BMT_Conference_Note_91158): normalize date, clean up evidence sentence
$ws.Range("E6").Value = "'2021-12-15"
$ws.Range("F6").Value = "['12/15: SPEP IgG kappa M spike of 0.1, Kappa light chain 1.24, lambda light chain 0.72. Kappa/ lambda ratio 1.72.']"

# Row 7 (now 2024-06-25_00:00:00.000_Progress_Notes_91596, kappa <0.06 mg/dL): normalize date
$ws.Range("E7").Value = "'2024-04-08"

# Row 8 (now 2024-04-17_00:00:00.000_Progress_Notes_91591): normalize date
$ws.Range("E8").Value = "'2024-02-22"

# These two cells stay empty both before and after the edit; make sure the
# row-delete/shift above didn't leave behind a stray empty-string value.
$ws.Range("C5").ClearContents() | Out-Null
$ws.Range("D7").ClearContents() | Out-Null
